$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.711.58'
$ws.Range("E2").Value = '  -1.33%  '
$ws.Range("D3").Value = '1.811.83'
$ws.Range("E3").Value = '  -1.44%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").Formula = "'229.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.18%  '
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Formula = "'39.51"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -9.10%  '
$ws.Range("E9").Value = '  +3.10%  '
$ws.Range("D10").Formula = "'0.0679"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.52%  '
$ws.Range("D11").Formula = "'0.0987"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.88%  '
$ws.Range("D12").Value = '2.072.88'
$ws.Range("E12").Value = '  -1.59%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Formula = "'11.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.35%  '
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Formula = "'0.663"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.46%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.794.34'
$ws.Range("E15").Value = '  -2.33%  '
$ws.Range("E16").Value = '  -2.98%  '
$ws.Range("D17").Value = '34.692.36'
$ws.Range("E17").Value = '  -1.85%  '
$ws.Range("D18").Formula = "'69.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.41%  '
$ws.Range("E19").Value = '  -1.56%  '
$ws.Range("D20").Formula = "'239.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.47%  '
$ws.Range("D21").Formula = "'11.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.66%  '
$ws.Range("D22").Formula = "'4.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.61%  '
$ws.Range("E23").Value = '  +0.36%  '
$ws.Range("D24").Formula = "'2.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").Formula = "'173.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.54%  '
$ws.Range("E26").Value = '  -1.62%  '
$ws.Range("E27").Value = '  +2.30%  '
$ws.Range("D28").Formula = "'17.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.03%  '
$ws.Range("E29").Value = '  -4.43%  '
$ws.Range("E30").Value = '  -0.33%  '
$ws.Range("D31").Formula = "'3.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.93%  '
$ws.Range("D32").Formula = "'0.0543"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.66%  '
$ws.Range("E33").Value = '  -3.50%  '
$ws.Range("E34").Value = '  +10.47%  '
$ws.Range("E35").Value = '  -0.58%  '
$ws.Range("D36").Formula = "'0.684"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.18%  '
$ws.Range("E37").Value = '  +5.82%  '
$ws.Range("D38").Formula = "'90.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.69%  '
$ws.Range("D39").Value = '1.331.78'
$ws.Range("E39").Value = '  +0.88%  '
$ws.Range("E40").Value = '  -1.56%  '
$ws.Range("E41").Value = '  -3.43%  '
$ws.Range("E42").Value = '  -1.35%  '
$ws.Range("E43").Value = '  -4.31%  '
$ws.Range("E44").Value = '  -6.68%  '
$ws.Range("E45").Value = '  -2.06%  '
$ws.Range("D46").Formula = "'0.0519"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.06%  '
$ws.Range("E47").Value = '  -1.74%  '
$ws.Range("D48").Value = '1.992.00'
$ws.Range("E48").Value = '  -1.28%  '
$ws.Range("E49").Value = '  +0.27%  '
$ws.Range("D50").Formula = "'0.0660"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.43%  '
$ws.Range("D51").Formula = "'96.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.10%  '
